# Auto-update draw results: append the new Pick 4 draw row for 2025-10-09.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 23
$ws.Cells.Item($row, 1).Value = "'2025-10-09"
$ws.Cells.Item($row, 1).Style = "Normal"
$ws.Cells.Item($row, 2).Value = "Pick 4"
$ws.Cells.Item($row, 3).Value = "'251009"
$ws.Cells.Item($row, 3).Style = "Normal"
$ws.Cells.Item($row, 4).Value = "6-7-5-6"
$ws.Cells.Item($row, 5).Value = "2025-10-09T21:39:04.253+04:00"
